$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "50 najbardziej ekscytujących nastolatków w światowym futbolu wg FourFourTwo"
$ws.Range("B31").Value = "Dwóch piłkarzy Blaugrany w rankingu"
$ws.Range("C31").Value = "http://fcbarca.com/108084-50-najbardziej-ekscytujacych-nastolatkow-w-swiatowym-futbolu-wg-fourfourtwo.html"
$ws.Range("D31").Value = 1
